# Updated cryptos list on Sat Jan  6 14:47:06 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Force text storage for numeric-looking price values (matches source data which is stored as text)
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D16", "D17", "D20", "D22", "D23", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D38", "D39", "D42", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "44.140.64"
$ws.Range("E2").Value = "  +0.93%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.244.85"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "307.33"
$ws.Range("E5").Value = "  -2.65%  "

# Row 6 - Solana
$ws.Range("D6").Value = "96.31"
$ws.Range("E6").Value = "  -2.54%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.02%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.11%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -0.91%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "35.08"
$ws.Range("E10").Value = "  -3.10%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  -0.65%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "7.26"
$ws.Range("E12").Value = "  -1.18%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.27%  "

# Row 14 - was WrappedliquidstakedEther2.0, now WrappedEther
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.378.53"
$ws.Range("E14").Value = "  +6.30%  "

# Row 15 - was WrappedEther, now WrappedliquidstakedEther2.0
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.590.22"
$ws.Range("E15").Value = "  +0.51%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.831"
$ws.Range("E16").Value = "  -0.93%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "13.63"

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "44.086.37"
$ws.Range("E18").Value = "  +0.97%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +0.82%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "12.34"
$ws.Range("E20").Value = "  -4.03%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.36%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "65.44"
$ws.Range("E22").Value = "  +1.03%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "237.14"
$ws.Range("E23").Value = "  +1.64%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  -3.04%  "

# Row 25 - ImmutableX
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "  -1.83%  "

# Row 27 - InjectiveProtocol
$ws.Range("D27").Value = "39.12"
$ws.Range("E27").Value = "  +6.09%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "9.94"
$ws.Range("E28").Value = "  -3.16%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +0.99%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").Value = "  +1.40%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "20.05"
$ws.Range("E31").Value = "  +0.73%  "

# Row 32 - Monero
$ws.Range("D32").Value = "152.15"
$ws.Range("E32").Value = "  -3.68%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.0805"
$ws.Range("E33").Value = "  -3.22%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "3.33"
$ws.Range("E34").Value = "  +4.06%  "

# Row 35 - WEMIXToken
$ws.Range("D35").Value = "2.62"
$ws.Range("E35").Value = "  -2.26%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -0.03%  "

# Row 37 - Stellar
$ws.Range("E37").Value = "  +3.19%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -6.00%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "15.00"
$ws.Range("E39").Value = "  -7.10%  "

# Row 40 - NEARProtocol
$ws.Range("E40").Value = "  -4.70%  "

# Row 41 - RenderToken
$ws.Range("E41").Value = "  -4.51%  "

# Row 42 - VeChain
$ws.Range("D42").Value = "0.0300"
$ws.Range("E42").Value = "  -2.89%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.13%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.729.03"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45 - BitcoinSV
$ws.Range("D45").Value = "84.47"
$ws.Range("E45").Value = "  +5.02%  "

# Row 46 - Algorand
$ws.Range("E46").Value = "  -2.04%  "

# Row 47 - Aave
$ws.Range("D47").Value = "100.45"
$ws.Range("E47").Value = "  -0.93%  "

# Row 48 - THORChain
$ws.Range("D48").Value = "4.90"
$ws.Range("E48").Value = "  -3.88%  "

# Row 49 - ordi
$ws.Range("D49").Value = "69.66"
$ws.Range("E49").Value = "  -5.33%  "

# Row 50 - FraxShare
$ws.Range("D50").Value = "8.12"
$ws.Range("E50").Value = "  +0.21%  "

# Row 51 - MultiversX
$ws.Range("D51").Value = "54.54"
$ws.Range("E51").Value = "  -3.49%  "
